# Insert 4 new rows of data (for Fecha 45212, the new week's prices) right
# before the existing row 545 block, pushing the rest of the table down.
# Row 544 (Fecha 44536, Segunda) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 545..548; Excel copies formatting (incl. the date
# style on column D) down from the row immediately above, same as a manual
# "Insert Copied Cells" / row-insert in the UI.
$ws.Rows("545:548").Insert()

# The static "stall" metadata (Mercado ID, Mercado, Region, Codreg, Tipo,
# Producto ID, Producto, Categoria ID, Categoria, Variedad, Unidad,
# Origen, Kg/unidad) is identical across this whole block, so copy it from
# the row that used to be 545 and now sits at 549 after the insert.
$cols = @("A","B","C","E","F","G","H","I","J","K","Q","R","T")
foreach ($col in $cols) {
    $srcVal = $ws.Range("$col" + "549").Value2
    for ($r = 545; $r -le 548; $r++) {
        $ws.Range("$col$r").Value = $srcVal
    }
}

# New block of 4 rows for Fecha 45212 (2023-10-?? serial date), one per
# Calidad grade: Especial, Primera, Segunda, Tercera.
$ws.Range("D545").Value = 45212
$ws.Range("L545").Value = "Especial"
$ws.Range("M545").Value = 97
$ws.Range("N545").Value = 11000
$ws.Range("O545").Value = 12000
$ws.Range("P545").Value = 11588
$ws.Range("S545").Value = 1655

$ws.Range("D546").Value = 45212
$ws.Range("L546").Value = "Primera"
$ws.Range("M546").Value = 75
$ws.Range("N546").Value = 9000
$ws.Range("O546").Value = 9000
$ws.Range("P546").Value = 9000
$ws.Range("S546").Value = 1286

$ws.Range("D547").Value = 45212
$ws.Range("L547").Value = "Segunda"
$ws.Range("M547").Value = 57
$ws.Range("N547").Value = 6000
$ws.Range("O547").Value = 6000
$ws.Range("P547").Value = 6000
$ws.Range("S547").Value = 857

$ws.Range("D548").Value = 45212
$ws.Range("L548").Value = "Tercera"
$ws.Range("M548").Value = 50
$ws.Range("N548").Value = 4000
$ws.Range("O548").Value = 4000
$ws.Range("P548").Value = 4000
$ws.Range("S548").Value = 571
